# Revert "Merging 0.1.8 w VitalSigns"
#
# Restores the pre-merge content of the CIBMTR SNOMED myelodysplastic
# diseases ValueSet workbook: the "Include" sheet's friendlier tab name,
# the earlier Version/Status/Date metadata values, the original duplicate
# "Contact" rows, and removes the "Jurisdiction" row that the merge had
# introduced (which shifts Description/Purpose/Copyright/Immutable back
# up by one row).

$wb = $excel.ActiveWorkbook

# --- workbook.xml: rename the second sheet tab ---
$wsInclude = $wb.Worksheets.Item(2)
$wsInclude.Name = "Include from SNOMED CT"

# --- Metadata sheet: restore the old property values ---
$ws = $wb.ActiveSheet   # "Metadata" is the active/first sheet

$ws.Range("B3").Value = "0.1.6"                         # Version: 0.1.8 -> 0.1.6
$ws.Range("B6").Value = "active"                         # Status: draft -> active
$ws.Range("B8").Value = "2023-05-05T10:50:04-05:00"      # Date

# The two "Contact" rows both revert to the pre-merge placeholder text
$ws.Range("B10").Value = "No display for ContactDetail"
$ws.Range("B11").Value = "No display for ContactDetail"

# Drop the "Jurisdiction" row entirely (row 12) -- Description/Purpose/
# Copyright/Immutable each move up one row as a result, and the sheet's
# used range shrinks from A1:B16 to A1:B15.
$ws.Rows.Item(12).Delete()

Write-Output "reverted Merging 0.1.8 w VitalSigns"
